{"js": "// Collapse the arzt.paragraphs IF/ELSE/END-IF block back down to a plain\n// FOR / INS / END-FOR block (matching kk.paragraphs / part3.paragraphs),\n// restoring the a4.docx template and updating the intro copy for that loop.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the `{{FOR p2 IN arzt.paragraphs}}` \u2026 `{{END-FOR p2}}` block by its\n// literal placeholder text so the edit is resilient to any surrounding\n// paragraphs shifting around.\nlet forIdx = -1;\nlet endForIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"{{FOR p2 IN arzt.paragraphs}}\") {\n    forIdx = i;\n  } else if (forIdx !== -1 && endForIdx === -1 && t === \"{{END-FOR p2}}\") {\n    endForIdx = i;\n    break;\n  }\n}\n\nif (forIdx === -1 || endForIdx === -1) {\n  throw new Error(\"Could not locate arzt.paragraphs FOR/END-FOR block\");\n}\n\n// Inside that block we expect:\n//   {{IF $p2 === 'Haftungsausschluss (vom Patienten zu unterzeichnen)'}}\n//   {{INS $p2}}      <- bold run\n//   {{ELSE}}\n//   {{INS $p2}}      <- plain run (kept)\n//   {{END-IF}}\n// Remove the IF/ELSE/END-IF scaffolding plus the bold duplicate, leaving a\n// single plain `{{INS $p2}}` paragraph between FOR and END-FOR.\nconst toDelete = [];\nfor (let i = forIdx + 1; i < endForIdx; i++) {\n  const t = items[i].text.trim();\n  if (t === \"{{IF $p2 === 'Haftungsausschluss (vom Patienten zu unterzeichnen)'}}\" ||\n      t === \"{{ELSE}}\" ||\n      t === \"{{END-IF}}\") {\n    toDelete.push(items[i]);\n  } else if (t === \"{{INS $p2}}\") {\n    items[i].font.load(\"bold\");\n  }\n}\nawait context.sync();\n\nfor (let i = forIdx + 1; i < endForIdx; i++) {\n  const t = items[i].text.trim();\n  if (t === \"{{INS $p2}}\" && items[i].font.bold === true) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Collapse the arzt.paragraphs IF/ELSE/END-IF block back down to a plain\n# FOR / INS / END-FOR block (matching kk.paragraphs / part3.paragraphs),\n# restoring the a4.docx template and updating the intro copy for that loop.\n$d = $word.ActiveDocument\n\n$total = $d.Paragraphs.Count\n\n# Locate the `{{FOR p2 IN arzt.paragraphs}}` \u2026 `{{END-FOR p2}}` block by its\n# literal placeholder text so the edit is resilient to any surrounding\n# paragraphs shifting around.\n$forIdx = -1\n$endForIdx = -1\nfor ($i = 1; $i -le $total; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"{{FOR p2 IN arzt.paragraphs}}\") {\n        $forIdx = $i\n    } elseif ($forIdx -ge 0 -and $endForIdx -eq -1 -and $t -eq \"{{END-FOR p2}}\") {\n        $endForIdx = $i\n    }\n}\n\nif ($forIdx -eq -1 -or $endForIdx -eq -1) {\n    throw \"Could not locate arzt.paragraphs FOR/END-FOR block\"\n}\n\n# Inside that block we expect:\n#   {{IF $p2 === 'Haftungsausschluss (vom Patienten zu unterzeichnen)'}}\n#   {{INS $p2}}      <- bold run\n#   {{ELSE}}\n#   {{INS $p2}}      <- plain run (kept)\n#   {{END-IF}}\n# Remove the IF/ELSE/END-IF scaffolding plus the bold duplicate, leaving a\n# single plain `{{INS $p2}}` paragraph between FOR and END-FOR.\n$toDelete = @()\nfor ($i = $forIdx + 1; $i -lt $endForIdx; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"{{IF `$p2 === 'Haftungsausschluss (vom Patienten zu unterzeichnen)'}}\" -or\n        $t -eq \"{{ELSE}}\" -or\n        $t -eq \"{{END-IF}}\") {\n        $toDelete += $i\n    } elseif ($t -eq \"{{INS `$p2}}\" -and $p.Range.Font.Bold) {\n        $toDelete += $i\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$toDelete = $toDelete | Sort-Object -Descending\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
